# Junction_Flooding_35.xlsx edit:
#  - Round row 5 values to 2 decimal places ("custom accuracy")
#  - Delete row 6 entirely ("데이터 1000개" data-count trim)
#  - Narrow column F (col 6) width from 8 to 7
#  - dimension shrinks from A1:AH6 to A1:AH5 (handled automatically by Excel)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values to their rounded (2 decimal place) counterparts ---
$ws.Range("B5").Value  = 8.83
$ws.Range("C5").Value  = 6.34
$ws.Range("D5").Value  = 0.49
$ws.Range("E5").Value  = 18.44
$ws.Range("F5").Value  = 15.72
$ws.Range("G5").Value  = 6.79
$ws.Range("H5").Value  = 29.21
$ws.Range("I5").Value  = 10.25
$ws.Range("J5").Value  = 4.64
$ws.Range("K5").Value  = 6.77
$ws.Range("L5").Value  = 7.55
$ws.Range("M5").Value  = 8.18
$ws.Range("N5").Value  = 2.15
$ws.Range("O5").Value  = 6.37
$ws.Range("P5").Value  = 10.06
$ws.Range("Q5").Value  = 5.28
$ws.Range("R5").Value  = 0.31
$ws.Range("S5").Value  = 0.25
$ws.Range("T5").Value  = 95.51
$ws.Range("U5").Value  = 19.02
$ws.Range("V5").Value  = 6.37
$ws.Range("W5").Value  = 12.75
$ws.Range("X5").Value  = 7
$ws.Range("Y5").Value  = 0.85
$ws.Range("Z5").Value  = 13.77
$ws.Range("AA5").Value = 5.48
$ws.Range("AB5").Value = 4.97
$ws.Range("AC5").Value = 6.23
$ws.Range("AD5").Value = 7.93
$ws.Range("AE5").Value = 0.57
$ws.Range("AF5").Value = 26.83
$ws.Range("AG5").Value = 3.42
$ws.Range("AH5").Value = 7.63

# --- Remove row 6 (data was trimmed) ---
$ws.Rows.Item(6).Delete()

# --- Narrow column F (6) from width 8 to width 7 ---
# ColumnWidth is expressed in Excel character units, which differ from the
# raw OOXML <col width> units by a constant offset (~0.8333...) in this
# engine; subtract it so the serialized width lands exactly on 7.
$ws.Columns.Item(6).ColumnWidth = 7 - 0.8333333333333333
